$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the refreshed cryptos snapshot. Column D prices and the
# Coin/Link/Volume cells are stored as plain text in the sheet
# (values such as "67.329.70" are not valid numbers, and values
# like "1.00" must keep their literal text form). Any cell whose
# new value would otherwise be auto-parsed by Excel as a number
# is forced to Text format first, then restored to the default
# (General) style so it matches the original, unstyled cells.

$ws.Range("D2").Value = '67.329.70'
$ws.Range("E2").Value = '  -1.04%  '
$ws.Range("D3").Value = '3.330.83'
$ws.Range("E3").Value = '  +2.02%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.20'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.83%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '183.17'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.81%  '
$ws.Range("B7").Value = 'USDC'
$ws.Range("C7").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("B8").Value = 'XRP'
$ws.Range("C8").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.603'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.48%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.128'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.62%  '
$ws.Range("E10").Value = '  +0.78%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.406'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.39%  '
$ws.Range("D12").Value = '3.909.49'
$ws.Range("E12").Value = '  +1.96%  '
$ws.Range("E13").Value = '  -0.79%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.28'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.32%  '
$ws.Range("D15").Value = '67.529.99'
$ws.Range("E15").Value = '  -0.77%  '
$ws.Range("E16").Value = '  -0.30%  '
$ws.Range("D17").Value = '3.327.53'
$ws.Range("E17").Value = '  +1.58%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '444.08'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +6.79%  '
$ws.Range("E19").Value = '  +2.24%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.66'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.83%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.69'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.09%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '74.01'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.88%  '
$ws.Range("E23").Value = '  -0.19%  '
$ws.Range("D24").Value = '3.467.38'
$ws.Range("E24").Value = '  +1.85%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.511'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.64%  '
$ws.Range("E26").Value = '  +1.23%  '
$ws.Range("E27").Value = '  +2.08%  '
$ws.Range("E28").Value = '  -3.58%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.994'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.94%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.97'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.32%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.87'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.31%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.32'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.47%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.78'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.04%  '
$ws.Range("E35").Value = '  -1.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.49'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.06%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '161.38'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.24%  '
$ws.Range("E38").Value = '  -2.00%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '27.03'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.75%  '
$ws.Range("D40").Value = '2.797.20'
$ws.Range("E40").Value = '  +6.03%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.790'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.60%  '
$ws.Range("E42").Value = '  +0.17%  '
$ws.Range("E43").Value = '  -1.15%  '
$ws.Range("E44").Value = '  -0.68%  '
$ws.Range("E45").Value = '  -0.13%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '24.56'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.43%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.36'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.68%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '324.30'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.76%  '
$ws.Range("E49").Value = '  -0.30%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.982'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.76%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '30.98'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.33%  '
